# Applies updated cryptocurrency market data (price & 1h volume change) to Sheet1,
# matching the "Updated cryptos list" GitHub Actions commit.
# Numeric-looking price strings are prefixed with a literal apostrophe so Excel
# stores/keeps them as text (matching the original inlineStr cell type) instead of
# auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "97.040.42"
$ws.Range("E2").Value = "  +4.42%  "
$ws.Range("D3").Value = "3.142.32"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'242.18"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'613.18"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +2.65%  "
$ws.Range("D8").Value = "'0.386"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.141.53"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "'0.786"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "'0.199"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "96.800.57"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").Value = "'0.0000241"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "'34.36"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'5.51"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "3.714.65"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "3.131.59"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'524.22"
$ws.Range("E19").Value = "  +19.25%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'3.52"
$ws.Range("E20").Value = "  -7.41%  "
$ws.Range("D21").Value = "'14.70"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").Value = "'5.71"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").Value = "'0.0000195"
$ws.Range("E23").Value = "  -4.72%  "
$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'90.23"
$ws.Range("E25").Value = "  +5.39%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'5.51"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").Value = "'11.60"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "3.295.86"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'0.239"
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("D31").Value = "'0.176"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("D32").Value = "'0.125"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.04"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'0.879"
$ws.Range("E34").Value = "  -12.35%  "
$ws.Range("D35").Value = "'26.75"
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("D37").Value = "'7.40"
$ws.Range("E37").Value = "  -8.66%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "'481.88"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'0.440"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").Value = "'1.23"
$ws.Range("D43").Value = "'3.57"
$ws.Range("E43").Value = "  -10.41%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'3.19"
$ws.Range("E45").Value = "  -5.05%  "
$ws.Range("D46").Value = "'159.82"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'1.94"
$ws.Range("E47").Value = "  +5.41%  "
$ws.Range("D48").Value = "'0.707"
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("D49").Value = "'4.50"
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("D50").Value = "'44.44"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0323"
$ws.Range("E51").Value = "  -0.34%  "
